{"js": "// 1. Delete the \"Meta description: ...\" paragraph that follows the H1 title.\nlet paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst metaDescPrefix = \"Meta description: Read our review of Fruits and Fire\";\nlet metaPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(metaDescPrefix) === 0) {\n    metaPara = paras.items[i];\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// 2. Re-query paragraphs (collection is stale after the delete) and find the\n//    final paragraph, which holds the old \"Prompt: ...\" image-generation text.\nparas = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst promptText = \"Prompt: Create a cartoon-style feature image for the game \\\"Fruits and Fire\\\" that showcases a happy Maya warrior with glasses. The image should be bright and colorful, featuring the warrior holding a bowl of fruits on one hand and a torch on the other. The background could depict a tropical jungle or a fiery ambiance to represent the game's name. Make sure the image highlights the game's retro-style and exotic touch, while still being eye-catching and fun.\";\nconst newMetaText = \"Read our review of Fruits and Fire, a simple slot game from Synot Games with expanding wilds and scatters. Play for free and enjoy the fruit theme.\";\n\nlet promptPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === promptText) {\n    promptPara = paras.items[i];\n    break;\n  }\n}\n\nif (promptPara) {\n  // 3. Insert a brand-new bold paragraph with the title text right before it.\n  const newPara = promptPara.insertParagraph(\n    \"Play Fruits and Fire for Free - A Simple Slot Game from Synot Games\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n\n  // Re-write that paragraph's OOXML directly so it ends up with exactly the\n  // structure the diff shows (leading empty run + single bold run), instead\n  // of inheriting direct character formatting (e.g. italics) from its\n  // neighbouring paragraph.\n  const titleOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +\n    '<w:t>Play Fruits and Fire for Free - A Simple Slot Game from Synot Games</w:t>' +\n    '</w:r></w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  newPara.insertOoxml(titleOoxml, Word.InsertLocation.replace);\n  await context.sync();\n\n  // 4. Replace the old \"Prompt: ...\" text with the meta-description copy,\n  //    keeping the paragraph's existing italic run formatting intact.\n  const searchResults = promptPara.search(promptText, { matchCase: true });\n  searchResults.load(\"text\");\n  await context.sync();\n  if (searchResults.items.length > 0) {\n    searchResults.items[0].insertText(newMetaText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Delete the \"Meta description: ...\" paragraph that follows the H1 title.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Meta description\")) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2. Locate the final paragraph, which holds the old \"Prompt: ...\" image\n#    generation text.\n$promptText = \"Prompt: Create a cartoon-style feature image for the game \"\"Fruits and Fire\"\" that showcases a happy Maya warrior with glasses. The image should be bright and colorful, featuring the warrior holding a bowl of fruits on one hand and a torch on the other. The background could depict a tropical jungle or a fiery ambiance to represent the game's name. Make sure the image highlights the game's retro-style and exotic touch, while still being eye-catching and fun.\"\n$newMetaText = \"Read our review of Fruits and Fire, a simple slot game from Synot Games with expanding wilds and scatters. Play for free and enjoy the fruit theme.\"\n\n$count = $d.Paragraphs.Count\n$promptPara = $d.Paragraphs.Item($count)\n\n# 3. Insert a brand-new paragraph right before it, then fill it in with raw\n#    OOXML so it ends up with exactly the structure the diff shows (leading\n#    empty run + single bold run) instead of inheriting direct character\n#    formatting (e.g. italics) from the neighbouring \"Prompt\" paragraph.\n$promptPara.Range.InsertParagraphBefore()\n$newCount = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($newCount - 1)\n\n$titleOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +\n    '<w:t>Play Fruits and Fire for Free - A Simple Slot Game from Synot Games</w:t>' +\n    '</w:r></w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n$newPara.Range.InsertXML($titleOoxml)\n\n# 4. Replace the old \"Prompt: ...\" text with the meta-description copy,\n#    keeping the paragraph's existing italic run formatting intact.\n$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$find = $promptPara.Range.Find\n$find.ClearFormatting()\n$find.Execute($promptText, $false, $false, $false, $false, $false, $true, 1, $false, $newMetaText, 2)\n"}
